$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "74÷8="; New = "42÷2=" },
    @{ Old = "26÷5="; New = "18÷2=" },
    @{ Old = "35÷2="; New = "21÷3=" },
    @{ Old = "98÷8="; New = "45÷9=" },
    @{ Old = "96÷6="; New = "86÷4=" },
    @{ Old = "69÷4="; New = "57÷3=" },
    @{ Old = "32÷2="; New = "57÷4=" },
    @{ Old = "29÷9="; New = "26÷9=" },
    @{ Old = "84÷7="; New = "90÷3=" },
    @{ Old = "50÷6="; New = "49÷4=" },
    @{ Old = "65÷2="; New = "60÷7=" },
    @{ Old = "57÷7="; New = "22÷2=" },
    @{ Old = "44÷4="; New = "11÷8=" },
    @{ Old = "68÷4="; New = "44÷5=" },
    @{ Old = "82÷8="; New = "48÷7=" },
    @{ Old = "41÷4="; New = "53÷4=" },
    @{ Old = "53÷6="; New = "53÷8=" },
    @{ Old = "70÷5="; New = "43÷2=" },
    @{ Old = "94÷2="; New = "21÷5=" },
    @{ Old = "89÷4="; New = "52÷8=" },
    @{ Old = "19÷4="; New = "53÷3=" },
    @{ Old = "77÷3="; New = "70÷6=" },
    @{ Old = "44÷2="; New = "45÷7=" },
    @{ Old = "23÷6="; New = "87÷3=" },
    @{ Old = "59÷3="; New = "96÷4=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
